# Update InsideBet Data: Automatizado
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - Villarreal
$ws.Range("C4").Value = 24
$ws.Range("D4").Value = 15
$ws.Range("G4").Value = 45
$ws.Range("I4").Value = 19
$ws.Range("J4").Value = 48
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = "L D W L W"

# Row 20 - Levante
$ws.Range("C20").Value = 24
$ws.Range("F20").Value = 14
$ws.Range("H20").Value = 41
$ws.Range("I20").Value = -15
$ws.Range("K20").Value = 0.75
$ws.Range("L20").Value = "W D L L L"
$ws.Range("M20").Value = 20575
